$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.747.29'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '3.155.89'
$ws.Range('E3').Value = '  +1.99%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.49'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '3.155.39'
$ws.Range('E8').Value = '  +2.07%  '
$ws.Range('E9').Value = '  +4.39%  '
$ws.Range('E10').Value = '  +5.93%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.14'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.35%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.506'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +7.10%  '
$ws.Range('E13').Value = '  +12.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.25'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +8.80%  '
$ws.Range('D15').Value = '3.677.01'
$ws.Range('E15').Value = '  +2.45%  '
$ws.Range('D16').Value = '64.868.34'
$ws.Range('E16').Value = '  +0.88%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.25'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +7.19%  '
$ws.Range('D18').Value = '3.156.84'
$ws.Range('E18').Value = '  +2.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '519.13'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.73%  '
$ws.Range('E20').Value = '  +0.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.738'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +9.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '15.23'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +7.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.86'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.35'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.96%  '
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.93'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.64%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.79'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +9.56%  '
$ws.Range('E29').Value = '  +6.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '28.02'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +6.06%  '
$ws.Range('E31').Value = '  +0.11%  '
$ws.Range('E32').Value = '  +8.04%  '
$ws.Range('E33').Value = '  +3.76%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.12'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +9.37%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.59'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '55.83'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '490.20'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.30%  '
$ws.Range('E38').Value = '  +6.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0425'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.99'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.18%  '
$ws.Range('D41').Value = '3.114.85'
$ws.Range('E41').Value = '  +4.79%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.68'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.36%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.120'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.65%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.297'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +13.52%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.46'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +15.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '29.31'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.46%  '
$ws.Range('D47').Value = '0.0₃0581'
$ws.Range('E47').Value = '  +12.75%  '
$ws.Range('E48').Value = '  -0.02%  '
$ws.Range('E49').Value = '  +3.26%  '
$ws.Range('E50').Value = '  +10.72%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '119.32'
$ws.Range('D51').Style = 'Normal'
